$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths per diff (C, D, F, H change; A, B, E, G unchanged)
# Note: Excel's ColumnWidth setter stores width + 5/6 in the OOXML "width"
# attribute (Maximum-Digit-Width rounding), so subtract 5/6 here to land on
# the exact target values (60, 50, 17, 34) in the saved file.
$offset = 5/6
$ws.Columns.Item(3).ColumnWidth = 60 - $offset
$ws.Columns.Item(4).ColumnWidth = 50 - $offset
$ws.Columns.Item(6).ColumnWidth = 17 - $offset
$ws.Columns.Item(8).ColumnWidth = 34 - $offset

# New data for rows 2-9 (opportunity id, link, title, country, premium, applicants, duration, organization)
$data = @(
    @("1326974", "https://aiesec.org/opportunity/global-talent/1326974", "[Impact Brazil] - Cooking SubSystems Intern", "Joinville - Pirabeiraba, Joinville - SC, Brasil", "No", "0 applicants", "6 - 18 Months", "Whirlpool Corporation"),
    @("1326970", "https://aiesec.org/opportunity/global-talent/1326970", "Head of Ecommerce", "Galle, Sri Lanka", "No", "0 applicants", "6 - 18 Months", "Tallentire House (Pvt) Ltd"),
    @("1326834", "https://aiesec.org/opportunity/global-talent/1326834", "Marketing Manager", "Jamshedpur, Jharkhand, India", "No", "0 applicants", "3 - 6 Months", "ToWaSo Pvt.Ltd."),
    @("1326411", "https://aiesec.org/opportunity/global-talent/1326411", "Marketing Trainee Pharmaceutical", "Panamá, Provincia de Panamá, Panamá", "No", "44 applicants", "6 - 18 Months", "NOVARTIS"),
    @("1325033", "https://aiesec.org/opportunity/global-talent/1325033", "Junior Full-Stack Developer – AI & Web Projects (EU ONLY)", "Brussels, Belgium", "No", "91 applicants", "6 - 18 Months", "Eureka Resource Mining"),
    @("1317126", "https://aiesec.org/opportunity/global-talent/1317126", "Social Media Manager", "Ahangama, Sri Lanka", "No", "14 applicants", "9 - 12 Weeks", "Surfing Wombats"),
    @("1316641", "https://aiesec.org/opportunity/global-talent/1316641", "Tour Executive - Spanish", "Colombo, Sri Lanka", "No", "7 applicants", "6 - 18 Months", "Aitken Spence Travels (Pvt) Ltd"),
    @("1315099", "https://aiesec.org/opportunity/global-talent/1315099", "Management Control Analyst", "Panamá, Provincia de Panamá, Panamá", "No", "148 applicants", "6 - 18 Months", "Skechers Latin America LLC")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        if ($col -eq 1) {
            # Opportunity ID column holds text-formatted numeric strings in the source
            $ws.Cells.Item($row, $col).Value = "'" + $rowData[$j]
        } else {
            $ws.Cells.Item($row, $col).Value = $rowData[$j]
        }
    }
}

# Remove old rows 10-15 (table now only spans through row 9)
$ws.Range("A10:H15").Delete()
